$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Marks")

# Set H6 value to 0 (Lab 7 Branch&Bound score for the student)
$ws.Range("H6").Value = 0

# Set H7 value to "Not done" (comment/feedback text for Lab 7)
$ws.Range("H7").Value = "Not done"

# Update selection / view state
$ws.Range("I7:I14").Select()
$excel.ActiveWindow.ScrollColumn = 4
